$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.505.17"
$ws.Range("E2").Value = "  +6.57%  "

$ws.Range("D3").Value = "2.471.93"
$ws.Range("E3").Value = "  +4.37%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "487.95"
$ws.Range("E5").Value = "  +6.94%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.36"
$ws.Range("E6").Value = "  +12.37%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.513"
$ws.Range("E8").Value = "  +7.52%  "

$ws.Range("D9").Value = "2.492.05"
$ws.Range("E9").Value = "  +4.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.76"
$ws.Range("E10").Value = "  +10.74%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0979"
$ws.Range("E11").Value = "  +4.59%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.332"
$ws.Range("E12").Value = "  +7.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.123"
$ws.Range("E13").Value = "  +2.01%  "

$ws.Range("D14").Value = "2.909.45"
$ws.Range("E14").Value = "  +4.61%  "

$ws.Range("D15").Value = "56.460.66"
$ws.Range("E15").Value = "  +6.23%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.14"
$ws.Range("E16").Value = "  +8.77%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000136"
$ws.Range("E17").Value = "  +6.11%  "

$ws.Range("D18").Value = "2.491.50"
$ws.Range("E18").Value = "  +4.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.54"
$ws.Range("E19").Value = "  +10.50%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.13"
$ws.Range("E20").Value = "  +10.02%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "318.27"
$ws.Range("E21").Value = "  +4.85%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.995"
$ws.Range("E22").Value = "  -0.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.81"
$ws.Range("E23").Value = "  +10.65%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "58.53"
$ws.Range("E24").Value = "  +5.74%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.410"
$ws.Range("E25").Value = "  +8.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  +10.19%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.88%  "

$ws.Range("D28").Value = "2.595.36"
$ws.Range("E28").Value = "  +4.36%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.59"
$ws.Range("E29").Value = "  +8.70%  "

$ws.Range("D30").Value = "0.0₃0784"
$ws.Range("E30").Value = "  +10.98%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.31%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "149.03"
$ws.Range("E32").Value = "  +2.98%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.18"
$ws.Range("E33").Value = "  +4.42%  "

$ws.Range("E34").Value = "  +7.11%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.20"
$ws.Range("E35").Value = "  +6.02%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.14"
$ws.Range("E36").Value = "  +9.83%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.72"
$ws.Range("E37").Value = "  +7.62%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.857"
$ws.Range("E38").Value = "  +9.76%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.14"
$ws.Range("E39").Value = "  +3.90%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.51"
$ws.Range("E40").Value = "  +9.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0559"
$ws.Range("E41").Value = "  +8.37%  "

$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.997"
$ws.Range("E42").Value = "  +0.15%  "

$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.611"
$ws.Range("E43").Value = "  +4.57%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.33"
$ws.Range("E44").Value = "  +9.86%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.80"
$ws.Range("E45").Value = "  +16.31%  "

$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0921"
$ws.Range("E46").Value = "  +7.65%  "

$ws.Range("B47").Value = "WhiteBITCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.21"
$ws.Range("E47").Value = "  +0.74%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0228"
$ws.Range("E48").Value = "  +7.00%  "

$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "256.31"
$ws.Range("E49").Value = "  +19.05%  "

$ws.Range("D50").Value = "1.894.93"
$ws.Range("E50").Value = "  -1.19%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.62"
$ws.Range("E51").Value = "  +8.40%  "

